# Replace the small "Greeting/Communication/ids" practice table with a
# "sr no / name / lastname / age" table (4 columns x 5 data rows), and
# remove the thin-border formatting that was applied to the old table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "sr no "
$ws.Range("B1").Value = "name "
$ws.Range("C1").Value = "lastname "
$ws.Range("D1").Value = "age "

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "lavanya "
$ws.Range("C2").Value = "ghadge"
$ws.Range("D2").Value = 14

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "roshan"
$ws.Range("C3").Value = "ghadge"
$ws.Range("D3").Value = 24

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "sanjay "
$ws.Range("C4").Value = "ghadge"
$ws.Range("D4").Value = 54

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "vijaya"
$ws.Range("C5").Value = "ghadge"
$ws.Range("D5").Value = 40

# The old A1:C4 block had a thin box border (style index 1); the new
# table carries no border formatting, so strip it back to the default.
$ws.Range("A1:C4").Borders.LineStyle = -4142

# Match the saved selection/extent of the edited table.
$ws.Range("A1:D5").Select()

Write-Output "done"
